$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "41.878.99"
Set-TextValue "E2" "  +5.62%  "
Set-TextValue "D3" "2.230.81"
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "231.70"
Set-TextValue "E5" "  +2.10%  "
Set-TextValue "E6" "  +0.31%  "
Set-TextValue "D7" "61.49"
Set-TextValue "E7" "  -2.71%  "
Set-TextValue "E8" "  +0.06%  "
Set-TextValue "D9" "0.401"
Set-TextValue "E9" "  +2.76%  "
Set-TextValue "D10" "59.16"
Set-TextValue "E10" "  +1.31%  "
Set-TextValue "D11" "0.0890"
Set-TextValue "E11" "  +4.30%  "
Set-TextValue "E12" "  -0.11%  "
Set-TextValue "D13" "2.561.09"
Set-TextValue "E13" "  +2.79%  "
Set-TextValue "D14" "15.65"
Set-TextValue "E14" "  -1.54%  "
Set-TextValue "D15" "22.03"
Set-TextValue "E15" "  +1.25%  "
Set-TextValue "D16" "0.800"
Set-TextValue "E16" "  -0.95%  "
Set-TextValue "D17" "5.59"
Set-TextValue "E17" "  +1.93%  "
Set-TextValue "D18" "2.248.90"
Set-TextValue "E18" "  +3.66%  "
Set-TextValue "D19" "41.785.29"
Set-TextValue "E19" "  +5.48%  "
Set-TextValue "D20" "72.08"
Set-TextValue "E20" "  +0.43%  "
Set-TextValue "D21" "0.0₃0894"
Set-TextValue "E21" "  -2.79%  "
Set-TextValue "D22" "6.05"
Set-TextValue "E22" "  +0.57%  "
Set-TextValue "D23" "251.26"
Set-TextValue "E23" "  +10.33%  "
Set-TextValue "E24" "  +0.00%  "
Set-TextValue "B25" "Toncoin"
Set-TextValue "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D25" "2.37"
Set-TextValue "E25" "  +3.38%  "
Set-TextValue "B26" "PancakeSwap"
Set-TextValue "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "2.39"
Set-TextValue "E26" "  +2.34%  "
Set-TextValue "D27" "9.70"
Set-TextValue "E27" "  +2.34%  "
Set-TextValue "D28" "0.142"
Set-TextValue "E28" "  +2.78%  "
Set-TextValue "D29" "167.18"
Set-TextValue "E29" "  -2.16%  "
Set-TextValue "D30" "19.94"
Set-TextValue "E30" "  +1.25%  "
Set-TextValue "E31" "  -1.74%  "
Set-TextValue "E32" "  -1.49%  "
Set-TextValue "E33" "  -0.12%  "
Set-TextValue "D34" "5.00"
Set-TextValue "E34" "  +6.42%  "
Set-TextValue "D35" "4.67"
Set-TextValue "E35" "  +3.56%  "
Set-TextValue "D36" "0.0635"
Set-TextValue "E36" "  +3.22%  "
Set-TextValue "D37" "6.65"
Set-TextValue "E37" "  -4.71%  "
Set-TextValue "E38" "  -3.62%  "
Set-TextValue "E39" "  -0.82%  "
Set-TextValue "D40" "0.000255"
Set-TextValue "E40" "  +30.55%  "
Set-TextValue "E41" "  -0.15%  "
Set-TextValue "D42" "0.0240"
Set-TextValue "E42" "  +5.21%  "
Set-TextValue "D43" "4.81"
Set-TextValue "E43" "  -1.73%  "
Set-TextValue "D44" "8.60"
Set-TextValue "E44" "  +8.87%  "
Set-TextValue "D45" "0.0985"
Set-TextValue "E45" "  +7.47%  "
Set-TextValue "D46" "1.23"
Set-TextValue "E46" "  +1.52%  "
Set-TextValue "D47" "98.92"
Set-TextValue "E47" "  -3.37%  "
Set-TextValue "D48" "1.477.52"
Set-TextValue "E48" "  -2.42%  "
Set-TextValue "B49" "HuobiToken"
Set-TextValue "C49" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D49" "2.81"
Set-TextValue "E49" "  +0.18%  "
Set-TextValue "B50" "InjectiveProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D50" "16.46"
Set-TextValue "E50" "  -6.91%  "
Set-TextValue "D51" "52.39"
Set-TextValue "E51" "  +8.49%  "
